# Apply the cryptos-list data refresh described by the commit diff.
# For column D (Price) values we force the cell to Text format before
# assigning, because several price strings (e.g. "563.84") would
# otherwise be auto-parsed by Excel as numbers, changing their stored
# representation (and losing formatting such as "62.247.72"-style
# thousand-grouped values). The NumberFormat is reset to the default
# ("Normal" style) right after the assignment so no visible formatting
# change is left behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '62.247.72'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.53%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.424.00'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.06%  '

# Row 4
$ws.Range("E4").Value = '  -0.07%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '563.84'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.21%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.46'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.40%  '

# Row 7
$ws.Range("E7").Value = '  +0.01%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.532'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.74%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.421.49'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.91%  '

# Row 10
$ws.Range("E10").Value = '  +1.87%  '

# Row 11
$ws.Range("E11").Value = '  -2.14%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.37'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.47%  '

# Row 13
$ws.Range("E13").Value = '  +0.75%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.97'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.93%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000177'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +5.28%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.861.93'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.02%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '62.079.52'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.20%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.424.67'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.89%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.35'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.41%  '

# Row 20
$ws.Range("E20").Value = '  +1.28%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '323.94'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.05%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.75'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.81%  '

# Row 23
$ws.Range("E23").Value = '  +0.01%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '65.58'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.93%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.71'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.71%  '

# Row 26
$ws.Range("E26").Value = '  +0.81%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '583.54'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +12.78%  '

# Row 28
$ws.Range("B28").Value = 'Binance-PegBSC-USD'
$ws.Range("C28").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.00'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.20%  '

# Row 29
$ws.Range("B29").Value = 'WrappedeETH'
$ws.Range("C29").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.528.47'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.43%  '

# Row 30
$ws.Range("B30").Value = 'PEPE'
$ws.Range("C30").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0945'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +5.08%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.46'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +5.10%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.27'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.94%  '

# Row 33
$ws.Range("E33").Value = '  +0.41%  '

# Row 34
$ws.Range("E34").Value = '  +2.13%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.56'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.87%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.73'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.93%  '

# Row 37
$ws.Range("E37").Value = '  +0.07%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.80'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.50%  '

# Row 39
$ws.Range("E39").Value = '  +1.56%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '153.58'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +4.56%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '18.68'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.02%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.83'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.03%  '

# Row 43
$ws.Range("E43").Value = '  -0.16%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.33'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +7.87%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '150.29'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.68%  '

# Row 46
$ws.Range("E46").Value = '  +1.62%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0538'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.65%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '20.37'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.77%  '

# Row 49
$ws.Range("E49").Value = '  +2.29%  '

# Row 50
$ws.Range("E50").Value = '  +2.12%  '

# Row 51
$ws.Range("E51").Value = '  +2.02%  '
